$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 4): give it a full (all four sides) medium border, same
#    as the rest of the table, by copying the formatting from a data cell
#    that already has that border (B5).
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B4:E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2. SHA3 row: "Tested?" already Yes, but "Complete?" flips from No to Yes.
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "Yes"

# ---------------------------------------------------------------------------
# 3. Add four new rows (22-25) describing new outstanding work items. Grab
#    the formatting of an existing row that already has the same visual
#    pattern (row 18) and then fill in the new text.
# ---------------------------------------------------------------------------
$ws.Range("B18:E18").Copy() | Out-Null
$ws.Range("B22:E25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B22").Value = "Write functions for Jump/JumpI"
$ws.Range("B23").Value = "Add functions for SelfDestruct/Stop"
$ws.Range("B24").Value = "Add functionality to identify function hashes"
$ws.Range("B25").Value = "Sort NumberFormatExceptions in arithmetic ops"

$ws.Range("D22").Value = "Not Started"
$ws.Range("D23").Value = "Not Started"
$ws.Range("D24").Value = "Not Started"
$ws.Range("D25").Value = "Not Started"

$ws.Range("C22").Value = "No"
$ws.Range("C23").Value = "No"
$ws.Range("C24").Value = "No"
$ws.Range("C25").Value = "No"

$ws.Range("E22").Value = "No"
$ws.Range("E23").Value = "No"
$ws.Range("E24").Value = "No"
$ws.Range("E25").Value = "No"

# ---------------------------------------------------------------------------
# 4. Extend the conditional formatting that colours the status columns so it
#    covers the four new rows as well.
# ---------------------------------------------------------------------------
$newRange = $ws.Range("C5:E25")
$conditions = $ws.Range("C5:E21").FormatConditions
for ($i = 1; $i -le $conditions.Count; $i++) {
    $conditions.Item($i).ModifyAppliesToRange($newRange)
}

# ---------------------------------------------------------------------------
# 5. Widen column B so the new, longer task descriptions fit, and move the
#    active selection to reflect where the user ended up working.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 38.166666666666664

$ws.Range("B12").Select() | Out-Null
